$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 462 (pushes old rows 462-503 down to 464-505)
$ws.Range("A462:A463").EntireRow.Insert()

# Update row 460 (existing row, new week's data)
$ws.Range("D460").Value = 44461
$ws.Range("M460").Value = 580
$ws.Range("P460").Value = 13517
$ws.Range("S460").Value = 1931

# Update row 461 (existing row, new week's data)
$ws.Range("D461").Value = 44461
$ws.Range("M461").Value = 650
$ws.Range("N461").Value = 11000
$ws.Range("P461").Value = 11462
$ws.Range("S461").Value = 1637

# Fill in newly inserted row 462 (this takes the values that old row 460 had: Especial)
$ws.Range("A462").Value = 9
$ws.Range("B462").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C462").Value = "Metropolitana"
$ws.Range("D462").Value = 44357
$ws.Range("E462").Value = 13
$ws.Range("F462").Value = "Fruta"
$ws.Range("G462").Value = 100101
$ws.Range("H462").Value = "Berries"
$ws.Range("I462").Value = 100112025
$ws.Range("J462").Value = "Frutilla"
$ws.Range("K462").Value = "Sin especificar"
$ws.Range("L462").Value = "Especial"
$ws.Range("M462").Value = 500
$ws.Range("N462").Value = 13000
$ws.Range("O462").Value = 14000
$ws.Range("P462").Value = 13500
$ws.Range("Q462").Value = "$/bandeja 7 kilos"
$ws.Range("R462").Value = "Provincia de Melipilla"
$ws.Range("S462").Value = 1929
$ws.Range("T462").Value = 7

# Fill in newly inserted row 463 (this takes the values that old row 461 had: Primera)
$ws.Range("A463").Value = 9
$ws.Range("B463").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C463").Value = "Metropolitana"
$ws.Range("D463").Value = 44357
$ws.Range("E463").Value = 13
$ws.Range("F463").Value = "Fruta"
$ws.Range("G463").Value = 100101
$ws.Range("H463").Value = "Berries"
$ws.Range("I463").Value = 100112025
$ws.Range("J463").Value = "Frutilla"
$ws.Range("K463").Value = "Sin especificar"
$ws.Range("L463").Value = "Primera"
$ws.Range("M463").Value = 400
$ws.Range("N463").Value = 10000
$ws.Range("O463").Value = 12000
$ws.Range("P463").Value = 11000
$ws.Range("Q463").Value = "$/bandeja 7 kilos"
$ws.Range("R463").Value = "Provincia de Melipilla"
$ws.Range("S463").Value = 1571
$ws.Range("T463").Value = 7

Write-Host "Done"
